$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 400
$ws.Range("I4").Value = 400
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -286
$ws.Range("N4").ClearContents() | Out-Null

$ws.Range("H132").Value = 5884959
$ws.Range("I132").Value = 6252159.5
$ws.Range("K132").Value = 18756478.5
$ws.Range("M132").Value = -18753948.5

$ws.Range("H133").Value = 16447.824
$ws.Range("J133").Value = 16447.824
$ws.Range("L133").Value = 16447.824
$ws.Range("N133").Value = -26567.824

$ws.Range("H136").Value = 17961.316
$ws.Range("J136").Value = 17961.316
$ws.Range("L136").Value = 17961.316
$ws.Range("N136").Value = -28161.316

$ws.Range("H137").Value = 2177516.8
$ws.Range("I137").Value = 4352778
$ws.Range("J137").Value = 2255.1738
$ws.Range("K137").Value = 13058334
$ws.Range("L137").Value = 6765.5214
$ws.Range("M137").Value = -13055784
$ws.Range("N137").Value = -11865.5214

$ws.Range("H138").Value = 4424.273
$ws.Range("I138").Value = 2731.258
$ws.Range("J138").Value = 5196.0884
$ws.Range("K138").Value = 8193.773999999999
$ws.Range("L138").Value = 15588.2652
$ws.Range("M138").Value = -3053.773999999999
$ws.Range("N138").Value = -25868.2652

$ws.Range("H139").Value = 20720.84
$ws.Range("J139").Value = 20720.84
$ws.Range("L139").Value = 20720.84
$ws.Range("N139").Value = -31000.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 17841.166
$ws.Range("J44").Value = 20109.4
$ws.Range("L44").Value = 20109.4
$ws.Range("N44").Value = -21085.4

$ws.Range("H88").Value = 2264.0715
$ws.Range("I88").Value = 1638
$ws.Range("J88").Value = 2611.889
$ws.Range("K88").Value = 1638
$ws.Range("L88").Value = 2611.889
$ws.Range("M88").Value = -1232
$ws.Range("N88").Value = -3423.889

$ws.Range("H91").Value = 2264.0715
$ws.Range("I91").Value = 1638
$ws.Range("J91").Value = 2611.889
$ws.Range("K91").Value = 1638
$ws.Range("L91").Value = 2611.889
$ws.Range("M91").Value = -234
$ws.Range("N91").Value = -5419.889

$ws.Range("H135").Value = 19667.98
$ws.Range("J135").Value = 19667.98
$ws.Range("L135").Value = 19667.98
$ws.Range("N135").Value = -29807.98

$ws.Range("H137").Value = 18653.875
$ws.Range("J137").Value = 18653.875
$ws.Range("L137").Value = 18653.875
$ws.Range("N137").Value = -28853.875

$ws.Range("H139").Value = 17116.031
$ws.Range("J139").Value = 17116.031
$ws.Range("L139").Value = 17116.031
$ws.Range("N139").Value = -27396.031

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 5000
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 5000
$ws.Range("M8").Value = -4860
$ws.Range("N8").ClearContents() | Out-Null

$ws.Range("H107").Value = 2686.8125
$ws.Range("I107").Value = 2298.625
$ws.Range("K107").Value = 2298.625
$ws.Range("M107").Value = -378.625

$ws.Range("H134").Value = 2690.125
$ws.Range("I134").Value = 2547.976
$ws.Range("K134").Value = 7643.928
$ws.Range("M134").Value = -5108.928

$ws.Range("H137").Value = 20075.639
$ws.Range("J137").Value = 19743.258
$ws.Range("L137").Value = 19743.258
$ws.Range("N137").Value = -29943.258

$ws.Range("H138").Value = 18111.889
$ws.Range("J138").Value = 18111.889
$ws.Range("L138").Value = 18111.889
$ws.Range("N138").Value = -28391.889

$ws.Range("H140").Value = 16612.824
$ws.Range("J140").Value = 16612.824
$ws.Range("L140").Value = 16612.824
$ws.Range("N140").Value = -26972.824

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 46251.5
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4888

$ws.Range("H31").Value = 2044604.1
$ws.Range("I31").Value = 3127307.8
$ws.Range("J31").Value = 6573.7646
$ws.Range("K31").Value = 3127307.8
$ws.Range("L31").Value = 6573.7646
$ws.Range("M31").Value = -3127012.8
$ws.Range("N31").Value = -7163.7646

$ws.Range("H34").Value = 2044604.1
$ws.Range("I34").Value = 3127307.8
$ws.Range("J34").Value = 6573.7646
$ws.Range("K34").Value = 3127307.8
$ws.Range("L34").Value = 6573.7646
$ws.Range("M34").Value = -3127105.8
$ws.Range("N34").Value = -6977.7646

$ws.Range("H39").Value = 35942.668
$ws.Range("I39").Value = 13000
$ws.Range("J39").Value = 47414
$ws.Range("K39").Value = 13000
$ws.Range("L39").Value = 47414
$ws.Range("M39").Value = -12609
$ws.Range("N39").Value = -48196

$ws.Range("H49").Value = 35942.668
$ws.Range("I49").Value = 13000
$ws.Range("J49").Value = 47414
$ws.Range("K49").Value = 13000
$ws.Range("L49").Value = 47414
$ws.Range("M49").Value = -12818
$ws.Range("N49").Value = -47778

$ws.Range("H58").Value = 10873176
$ws.Range("I58").Value = 1800.25
$ws.Range("J58").Value = 22732858
$ws.Range("K58").Value = 1800.25
$ws.Range("L58").Value = 22732858
$ws.Range("M58").Value = -1597.25
$ws.Range("N58").Value = -22733264

$ws.Range("H136").Value = 10873176
$ws.Range("I136").Value = 1800.25
$ws.Range("J136").Value = 22732858
$ws.Range("K136").Value = 5400.75
$ws.Range("L136").Value = 68198574
$ws.Range("M136").Value = -2850.75
$ws.Range("N136").Value = -68203674

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 928.625
$ws.Range("I18").Value = 486
$ws.Range("J18").Value = 1666.3334
$ws.Range("K18").Value = 1458
$ws.Range("L18").Value = 4999.0002
$ws.Range("M18").Value = -1289
$ws.Range("N18").Value = -5337.0002

$ws.Range("H56").Value = 5446.273
$ws.Range("I56").Value = 5446.273
$ws.Range("K56").Value = 5446.273
$ws.Range("M56").Value = -4916.273

$ws.Range("H107").Value = 1514.4286
$ws.Range("I107").Value = 233
$ws.Range("J107").Value = 2475.5
$ws.Range("K107").Value = 699
$ws.Range("L107").Value = 7426.5
$ws.Range("M107").Value = 1221
$ws.Range("N107").Value = -11266.5

$ws.Range("H131").Value = 1446.9623
$ws.Range("I131").Value = 7466.6665
$ws.Range("J131").Value = 1085.78
$ws.Range("K131").Value = 22399.9995
$ws.Range("L131").Value = 3257.34
$ws.Range("M131").Value = -17359.9995
$ws.Range("N131").Value = -13337.34

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4099.4346
$ws.Range("I132").Value = 4469.731
$ws.Range("J132").Value = 3618.05
$ws.Range("K132").Value = 13409.193
$ws.Range("L132").Value = 10854.15
$ws.Range("M132").Value = -10879.193
$ws.Range("N132").Value = -15914.15

$ws.Range("H133").Value = 19472.086
$ws.Range("J133").Value = 19472.086
$ws.Range("L133").Value = 19472.086
$ws.Range("N133").Value = -29592.086

$ws.Range("H135").Value = 18309.025
$ws.Range("J135").Value = 18309.025
$ws.Range("L135").Value = 18309.025
$ws.Range("N135").Value = -28449.025

$ws.Range("H140").Value = 17100.857
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 17100.857
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 17100.857
$ws.Range("M140").ClearContents() | Out-Null
$ws.Range("N140").Value = -27460.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2314.5
$ws.Range("I7").Value = 1471.75
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 1471.75
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -1359.75
$ws.Range("N7").Value = -4224

$ws.Range("H55").Value = 1400.4546
$ws.Range("I55").Value = 288.125
$ws.Range("J55").Value = 4366.6665
$ws.Range("K55").Value = 288.125
$ws.Range("L55").Value = 4366.6665
$ws.Range("M55").Value = -115.125
$ws.Range("N55").Value = -4712.6665

$ws.Range("H106").Value = 26754
$ws.Range("J106").Value = 26754
$ws.Range("L106").Value = 26754
$ws.Range("N106").Value = -29278

$ws.Range("H126").Value = 2314.5
$ws.Range("I126").Value = 1471.75
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 4415.25
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -1945.25
$ws.Range("N126").Value = -16940

$ws.Range("H132").Value = 3113.5862
$ws.Range("I132").Value = 2839.3
$ws.Range("J132").Value = 3257.9473
$ws.Range("K132").Value = 8517.900000000001
$ws.Range("L132").Value = 9773.841899999999
$ws.Range("M132").Value = -5987.900000000001
$ws.Range("N132").Value = -14833.8419

$ws.Range("H135").Value = 20652.258
$ws.Range("J135").Value = 20652.258
$ws.Range("L135").Value = 20652.258
$ws.Range("N135").Value = -30792.258

$ws.Range("H139").Value = 22375.406
$ws.Range("J139").Value = 22375.406
$ws.Range("L139").Value = 22375.406
$ws.Range("N139").Value = -32655.406

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("N104").ClearContents() | Out-Null

$ws.Range("H126").Value = 47977.305
$ws.Range("I126").Value = 69839.39999999999
$ws.Range("J126").Value = 6985.875
$ws.Range("K126").Value = 209518.2
$ws.Range("L126").Value = 20957.625
$ws.Range("M126").Value = -207048.2
$ws.Range("N126").Value = -25897.625

$ws.Range("H132").Value = 1934118.1
$ws.Range("I132").Value = 2635239.2
$ws.Range("J132").Value = 31075.572
$ws.Range("K132").Value = 7905717.600000001
$ws.Range("L132").Value = 93226.716
$ws.Range("M132").Value = -7903187.600000001
$ws.Range("N132").Value = -98286.716

$ws.Range("H135").Value = 20059.797
$ws.Range("J135").Value = 20059.797
$ws.Range("L135").Value = 20059.797
$ws.Range("N135").Value = -30199.797
